$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-28 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-03-01 Sunday", 2) | Out-Null
$d.Content.Find.Execute("24+32=56", $true, $false, $false, $false, $false, $true, 1, $false, "40+21=61", 2) | Out-Null
$d.Content.Find.Execute("40+47=87", $true, $false, $false, $false, $false, $true, 1, $false, "24+70=94", 2) | Out-Null
$d.Content.Find.Execute("5+25=30", $true, $false, $false, $false, $false, $true, 1, $false, "48-20=28", 2) | Out-Null
$d.Content.Find.Execute("53-47=6", $true, $false, $false, $false, $false, $true, 1, $false, "64+7=71", 2) | Out-Null
$d.Content.Find.Execute("66+20=86", $true, $false, $false, $false, $false, $true, 1, $false, "12+7=19", 2) | Out-Null
$d.Content.Find.Execute("87+9=96", $true, $false, $false, $false, $false, $true, 1, $false, "64-25=39", 2) | Out-Null
$d.Content.Find.Execute("88+1=89", $true, $false, $false, $false, $false, $true, 1, $false, "98-38=60", 2) | Out-Null
$d.Content.Find.Execute("5+35=40", $true, $false, $false, $false, $false, $true, 1, $false, "45+6=51", 2) | Out-Null
$d.Content.Find.Execute("33-1=32", $true, $false, $false, $false, $false, $true, 1, $false, "87-9=78", 2) | Out-Null
$d.Content.Find.Execute("78-15=63", $true, $false, $false, $false, $false, $true, 1, $false, "36+49=85", 2) | Out-Null
$d.Content.Find.Execute("81+16=97", $true, $false, $false, $false, $false, $true, 1, $false, "1+31=32", 2) | Out-Null
$d.Content.Find.Execute("74-67=7", $true, $false, $false, $false, $false, $true, 1, $false, "92-77=15", 2) | Out-Null
$d.Content.Find.Execute("46-38=8", $true, $false, $false, $false, $false, $true, 1, $false, "73-33=40", 2) | Out-Null
$d.Content.Find.Execute("18+24=42", $true, $false, $false, $false, $false, $true, 1, $false, "44+34=78", 2) | Out-Null
$d.Content.Find.Execute("38+44=82", $true, $false, $false, $false, $false, $true, 1, $false, "67+26=93", 2) | Out-Null
$d.Content.Find.Execute("51+20=71", $true, $false, $false, $false, $false, $true, 1, $false, "18+42=60", 2) | Out-Null
$d.Content.Find.Execute("49-44=5", $true, $false, $false, $false, $false, $true, 1, $false, "87+0=87", 2) | Out-Null
$d.Content.Find.Execute("40+18=58", $true, $false, $false, $false, $false, $true, 1, $false, "15+11=26", 2) | Out-Null
$d.Content.Find.Execute("90-60=30", $true, $false, $false, $false, $false, $true, 1, $false, "13+2=15", 2) | Out-Null
$d.Content.Find.Execute("80-11=69", $true, $false, $false, $false, $false, $true, 1, $false, "22-0=22", 2) | Out-Null
$d.Content.Find.Execute("28-14=14", $true, $false, $false, $false, $false, $true, 1, $false, "41-28=13", 2) | Out-Null
$d.Content.Find.Execute("93-25=68", $true, $false, $false, $false, $false, $true, 1, $false, "45+27=72", 2) | Out-Null
$d.Content.Find.Execute("96-45=51", $true, $false, $false, $false, $false, $true, 1, $false, "85-55=30", 2) | Out-Null
$d.Content.Find.Execute("55+8=63", $true, $false, $false, $false, $false, $true, 1, $false, "23+13=36", 2) | Out-Null
$d.Content.Find.Execute("38+50=88", $true, $false, $false, $false, $false, $true, 1, $false, "64+14=78", 2) | Out-Null
$d.Content.Find.Execute("76-42=34", $true, $false, $false, $false, $false, $true, 1, $false, "16-16=0", 2) | Out-Null
$d.Content.Find.Execute("7+13=20", $true, $false, $false, $false, $false, $true, 1, $false, "35-18=17", 2) | Out-Null
$d.Content.Find.Execute("32-28=4", $true, $false, $false, $false, $false, $true, 1, $false, "67-50=17", 2) | Out-Null
$d.Content.Find.Execute("12+13=25", $true, $false, $false, $false, $false, $true, 1, $false, "27+68=95", 2) | Out-Null
$d.Content.Find.Execute("39+46=85", $true, $false, $false, $false, $false, $true, 1, $false, "82-40=42", 2) | Out-Null
$d.Content.Find.Execute("37+53=90", $true, $false, $false, $false, $false, $true, 1, $false, "92-86=6", 2) | Out-Null
$d.Content.Find.Execute("15+83=98", $true, $false, $false, $false, $false, $true, 1, $false, "72-45=27", 2) | Out-Null
$d.Content.Find.Execute("11+29=40", $true, $false, $false, $false, $false, $true, 1, $false, "32+25=57", 2) | Out-Null
$d.Content.Find.Execute("76+17=93", $true, $false, $false, $false, $false, $true, 1, $false, "99-50=49", 2) | Out-Null
$d.Content.Find.Execute("82+8=90", $true, $false, $false, $false, $false, $true, 1, $false, "46+26=72", 2) | Out-Null
$d.Content.Find.Execute("58-24=34", $true, $false, $false, $false, $false, $true, 1, $false, "48+17=65", 2) | Out-Null
$d.Content.Find.Execute("34+13=47", $true, $false, $false, $false, $false, $true, 1, $false, "22+40=62", 2) | Out-Null
$d.Content.Find.Execute("76-44=32", $true, $false, $false, $false, $false, $true, 1, $false, "84-42=42", 2) | Out-Null
$d.Content.Find.Execute("72-41=31", $true, $false, $false, $false, $false, $true, 1, $false, "55-9=46", 2) | Out-Null
$d.Content.Find.Execute("90-80=10", $true, $false, $false, $false, $false, $true, 1, $false, "41-0=41", 2) | Out-Null
$d.Content.Find.Execute("61-20=41", $true, $false, $false, $false, $false, $true, 1, $false, "44+26=70", 2) | Out-Null
$d.Content.Find.Execute("24+63=87", $true, $false, $false, $false, $false, $true, 1, $false, "44-10=34", 2) | Out-Null
$d.Content.Find.Execute("4+47=51", $true, $false, $false, $false, $false, $true, 1, $false, "78+3=81", 2) | Out-Null
$d.Content.Find.Execute("11+15=26", $true, $false, $false, $false, $false, $true, 1, $false, "62-41=21", 2) | Out-Null
$d.Content.Find.Execute("74+14=88", $true, $false, $false, $false, $false, $true, 1, $false, "80-9=71", 2) | Out-Null
$d.Content.Find.Execute("47-12=35", $true, $false, $false, $false, $false, $true, 1, $false, "10+71=81", 2) | Out-Null
$d.Content.Find.Execute("23-2=21", $true, $false, $false, $false, $false, $true, 1, $false, "49+9=58", 2) | Out-Null
$d.Content.Find.Execute("52-48=4", $true, $false, $false, $false, $false, $true, 1, $false, "0+52=52", 2) | Out-Null
$d.Content.Find.Execute("5+28=33", $true, $false, $false, $false, $false, $true, 1, $false, "97-91=6", 2) | Out-Null
$d.Content.Find.Execute("93-87=6", $true, $false, $false, $false, $false, $true, 1, $false, "45+41=86", 2) | Out-Null
$d.Content.Find.Execute("11+31=42", $true, $false, $false, $false, $false, $true, 1, $false, "67+2=69", 2) | Out-Null
$d.Content.Find.Execute("93-34=59", $true, $false, $false, $false, $false, $true, 1, $false, "64-12=52", 2) | Out-Null
$d.Content.Find.Execute("6+63=69", $true, $false, $false, $false, $false, $true, 1, $false, "92-78=14", 2) | Out-Null
$d.Content.Find.Execute("11+69=80", $true, $false, $false, $false, $false, $true, 1, $false, "48-34=14", 2) | Out-Null
$d.Content.Find.Execute("1+14=15", $true, $false, $false, $false, $false, $true, 1, $false, "36+51=87", 2) | Out-Null
$d.Content.Find.Execute("14+25=39", $true, $false, $false, $false, $false, $true, 1, $false, "67-6=61", 2) | Out-Null
$d.Content.Find.Execute("17+79=96", $true, $false, $false, $false, $false, $true, 1, $false, "61-35=26", 2) | Out-Null
$d.Content.Find.Execute("39-20=19", $true, $false, $false, $false, $false, $true, 1, $false, "42+19=61", 2) | Out-Null
$d.Content.Find.Execute("73+6=79", $true, $false, $false, $false, $false, $true, 1, $false, "80+2=82", 2) | Out-Null
$d.Content.Find.Execute("32-13=19", $true, $false, $false, $false, $false, $true, 1, $false, "99-35=64", 2) | Out-Null
$d.Content.Find.Execute("40-9=31", $true, $false, $false, $false, $false, $true, 1, $false, "40+11=51", 2) | Out-Null
$d.Content.Find.Execute("30+40=70", $true, $false, $false, $false, $false, $true, 1, $false, "52-22=30", 2) | Out-Null
$d.Content.Find.Execute("74-48=26", $true, $false, $false, $false, $false, $true, 1, $false, "30+16=46", 2) | Out-Null
$d.Content.Find.Execute("35+34=69", $true, $false, $false, $false, $false, $true, 1, $false, "82-80=2", 2) | Out-Null
$d.Content.Find.Execute("57-26=31", $true, $false, $false, $false, $false, $true, 1, $false, "25+43=68", 2) | Out-Null
$d.Content.Find.Execute("26+5=31", $true, $false, $false, $false, $false, $true, 1, $false, "46-41=5", 2) | Out-Null
$d.Content.Find.Execute("88+2=90", $true, $false, $false, $false, $false, $true, 1, $false, "50-40=10", 2) | Out-Null
$d.Content.Find.Execute("12+41=53", $true, $false, $false, $false, $false, $true, 1, $false, "33-4=29", 2) | Out-Null
$d.Content.Find.Execute("79-12=67", $true, $false, $false, $false, $false, $true, 1, $false, "44+7=51", 2) | Out-Null
$d.Content.Find.Execute("82-26=56", $true, $false, $false, $false, $false, $true, 1, $false, "47-9=38", 2) | Out-Null
$d.Content.Find.Execute("34-31=3", $true, $false, $false, $false, $false, $true, 1, $false, "73+14=87", 2) | Out-Null
$d.Content.Find.Execute("7+65=72", $true, $false, $false, $false, $false, $true, 1, $false, "44-32=12", 2) | Out-Null
$d.Content.Find.Execute("49+4=53", $true, $false, $false, $false, $false, $true, 1, $false, "79-29=50", 2) | Out-Null
$d.Content.Find.Execute("10+63=73", $true, $false, $false, $false, $false, $true, 1, $false, "17+66=83", 2) | Out-Null
$d.Content.Find.Execute("27+52=79", $true, $false, $false, $false, $false, $true, 1, $false, "75+1=76", 2) | Out-Null
$d.Content.Find.Execute("20-17=3", $true, $false, $false, $false, $false, $true, 1, $false, "70-5=65", 2) | Out-Null
$d.Content.Find.Execute("25+70=95", $true, $false, $false, $false, $false, $true, 1, $false, "19+1=20", 2) | Out-Null
$d.Content.Find.Execute("61-6=55", $true, $false, $false, $false, $false, $true, 1, $false, "20-8=12", 2) | Out-Null
$d.Content.Find.Execute("18+2=20", $true, $false, $false, $false, $false, $true, 1, $false, "29+21=50", 2) | Out-Null
$d.Content.Find.Execute("37+55=92", $true, $false, $false, $false, $false, $true, 1, $false, "96-85=11", 2) | Out-Null
$d.Content.Find.Execute("52-10=42", $true, $false, $false, $false, $false, $true, 1, $false, "88-1=87", 2) | Out-Null
$d.Content.Find.Execute("72+5=77", $true, $false, $false, $false, $false, $true, 1, $false, "52-3=49", 2) | Out-Null
$d.Content.Find.Execute("56+0=56", $true, $false, $false, $false, $false, $true, 1, $false, "43+51=94", 2) | Out-Null
$d.Content.Find.Execute("51-31=20", $true, $false, $false, $false, $false, $true, 1, $false, "95-15=80", 2) | Out-Null
$d.Content.Find.Execute("85-60=25", $true, $false, $false, $false, $false, $true, 1, $false, "99-40=59", 2) | Out-Null
$d.Content.Find.Execute("70-58=12", $true, $false, $false, $false, $false, $true, 1, $false, "85-43=42", 2) | Out-Null
$d.Content.Find.Execute("5+37=42", $true, $false, $false, $false, $false, $true, 1, $false, "29+19=48", 2) | Out-Null
$d.Content.Find.Execute("89-25=64", $true, $false, $false, $false, $false, $true, 1, $false, "17-12=5", 2) | Out-Null
$d.Content.Find.Execute("38+30=68", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=18", 2) | Out-Null
$d.Content.Find.Execute("81-62=19", $true, $false, $false, $false, $false, $true, 1, $false, "86+7=93", 2) | Out-Null
$d.Content.Find.Execute("85-4=81", $true, $false, $false, $false, $false, $true, 1, $false, "93-75=18", 2) | Out-Null
$d.Content.Find.Execute("17+17=34", $true, $false, $false, $false, $false, $true, 1, $false, "13+48=61", 2) | Out-Null
$d.Content.Find.Execute("32-1=31", $true, $false, $false, $false, $false, $true, 1, $false, "8+2=10", 2) | Out-Null
$d.Content.Find.Execute("17-1=16", $true, $false, $false, $false, $false, $true, 1, $false, "82-31=51", 2) | Out-Null
$d.Content.Find.Execute("80-7=73", $true, $false, $false, $false, $false, $true, 1, $false, "34+51=85", 2) | Out-Null
$d.Content.Find.Execute("23-4=19", $true, $false, $false, $false, $false, $true, 1, $false, "74-61=13", 2) | Out-Null
$d.Content.Find.Execute("99-46=53", $true, $false, $false, $false, $false, $true, 1, $false, "1+47=48", 2) | Out-Null
$d.Content.Find.Execute("30+32=62", $true, $false, $false, $false, $false, $true, 1, $false, "37+2=39", 2) | Out-Null
$d.Content.Find.Execute("62+14=76", $true, $false, $false, $false, $false, $true, 1, $false, "28-22=6", 2) | Out-Null
$d.Content.Find.Execute("56-35=21", $true, $false, $false, $false, $false, $true, 1, $false, "34+32=66", 2) | Out-Null
